$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, matching the header style used by the
# other header cells (B1:G1) - copy formats only (xlPasteFormats = -4122)
# from G1 so we reuse the existing style index instead of minting a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data rows for the new Save column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
